$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-69 down to 17-70.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with its data (matching the style/format of sibling rows).
$ws.Cells.Item(16, 1).Value2 = 8
$ws.Cells.Item(16, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value2 = "Coquimbo"
$ws.Cells.Item(16, 4).Value2 = 44607
$ws.Cells.Item(16, 5).Value2 = 4
$ws.Cells.Item(16, 6).Value2 = 100112030
$ws.Cells.Item(16, 7).Value2 = "Poroto granado"
$ws.Cells.Item(16, 8).Value2 = "Sin especificar"
$ws.Cells.Item(16, 9).Value2 = "Primera"
$ws.Cells.Item(16, 10).Value2 = 600
$ws.Cells.Item(16, 11).Value2 = 30000
$ws.Cells.Item(16, 12).Value2 = 31000
$ws.Cells.Item(16, 13).Value2 = 30500
$ws.Cells.Item(16, 14).Value2 = "`$/malla 25 kilos"
$ws.Cells.Item(16, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(16, 16).Value2 = 1220
$ws.Cells.Item(16, 17).Value2 = 25
$ws.Cells.Item(16, 18).Value2 = "Hortaliza"

# Apply the same date number format as the other rows in column D.
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
